$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.443.84'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.870.53'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('D5').Value = '''246.37'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '''0.4739'
$ws.Range('D8').Value = '''0.2918'
$ws.Range('E8').Value = '  +2.43%  '
$ws.Range('D9').Value = '''0.06498'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = '''22.03'
$ws.Range('E10').Value = '  +6.08%  '
$ws.Range('D11').Value = '''0.07712'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').Value = '''0.7389'
$ws.Range('E13').Value = '  +8.25%  '
$ws.Range('D14').Value = '1.871.12'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').Value = '''5.125'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '''274.42'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = '30.421.15'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '''1.000'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '2.118.35'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '''1.0000'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '''5.230'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').Value = '''6.172'
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').Value = '''9.308'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '''164.16'
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('D27').Value = '''18.83'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').Value = '''1.931'
$ws.Range('E28').Value = '  +1.99%  '
$ws.Range('D29').Value = '''0.09996'
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('D30').Value = '''1.367'
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').Value = '''1.502'
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').Value = '''4.306'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').Value = '''4.145'
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').Value = '''0.04833'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').Value = '''0.6984'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '''2.712'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01859'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''2.739'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''6.313'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''72.84'
$ws.Range('E41').Value = '  +3.85%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''1.970'
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '''0.4195'
$ws.Range('E43').Value = '  +3.32%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = '''0.9998'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '''0.8350'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '''102.07'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''9.224'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '''7.014'
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''932.21'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''35.37'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05642'
$ws.Range('E51').Value = '  +1.69%  '
